$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.75
$ws.Range("G2").Value = 1.87
$ws.Range("I2").Value = 5.9
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 1.41
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.45
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 1.84
$ws.Range("Q2").Value = 1.96
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.45
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 1.95
$ws.Range("W2").Value = 2.14
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 20
$ws.Range("Z2").Value = 980
$ws.Range("AB2").Value = 9.6
$ws.Range("AD2").Value = 22
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 80
$ws.Range("AJ2").Value = 22
$ws.Range("AK2").Value = 22
$ws.Range("AL2").Value = 48
$ws.Range("AM2").Value = 150
$ws.Range("AO2").Value = 110
$ws.Range("F3").Value = 2.1
$ws.Range("G3").Value = 2.72
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4.4
$ws.Range("J3").Value = 2.84
$ws.Range("K3").Value = 3.9
$ws.Range("L3").Value = 1.37
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 1.71
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.27
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.76
$ws.Range("U3").Value = 1.94
$ws.Range("V3").Value = 1.3
$ws.Range("W3").Value = 1.58
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 28
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AD3").Value = 18.5
$ws.Range("AE3").Value = 55
$ws.Range("AF3").Value = 17.5
$ws.Range("AG3").Value = 13.5
$ws.Range("AH3").Value = 21
$ws.Range("AJ3").Value = 38
$ws.Range("AK3").Value = 32
$ws.Range("AN3").Value = 24
$ws.Range("F4").Value = 1.74
$ws.Range("G4").Value = 1.83
$ws.Range("H4").Value = 4.8
$ws.Range("I4").Value = 5.3
$ws.Range("K4").Value = 4.7
$ws.Range("P4").Value = 2.12
$ws.Range("S4").Value = 2.84
$ws.Range("T4").Value = 1.72
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.23
$ws.Range("W4").Value = 2.2
$ws.Range("AI4").Value = 65
$ws.Range("G5").Value = 2.42
$ws.Range("H5").Value = 2.58
$ws.Range("I5").Value = 3.25
$ws.Range("J5").Value = 3.55
$ws.Range("N5").Value = 6.8
$ws.Range("O5").Value = 1.12
$ws.Range("Q5").Value = 1.37
$ws.Range("R5").Value = 1.84
$ws.Range("X5").Value = 46
$ws.Range("Y5").Value = 29
$ws.Range("F6").Value = 2.38
$ws.Range("G6").Value = 2.4
$ws.Range("I6").Value = 3.55
$ws.Range("J6").Value = 3.25
$ws.Range("K6").Value = 3.6
$ws.Range("N6").Value = 3.15
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 2.04
$ws.Range("R6").Value = 1.28
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.39
$ws.Range("W6").Value = 1.71
$ws.Range("Y6").Value = 12.5
$ws.Range("Z6").Value = 24
$ws.Range("AA6").Value = 65
$ws.Range("AB6").Value = 11
$ws.Range("AC6").Value = 7.8
$ws.Range("AD6").Value = 17
$ws.Range("AE6").Value = 46
$ws.Range("AF6").Value = 17
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 19
$ws.Range("AI6").Value = 65
$ws.Range("AJ6").Value = 38
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 55
$ws.Range("AM6").Value = 130
$ws.Range("AN6").Value = 30
$ws.Range("AO6").Value = 46
$ws.Range("F7").Value = 2.58
$ws.Range("G7").Value = 2.66
$ws.Range("H7").Value = 3.15
$ws.Range("I7").Value = 3.25
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 3.3
$ws.Range("L7").Value = 1.49
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 3.05
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 1.68
$ws.Range("Q7").Value = 2.3
$ws.Range("R7").Value = 1.25
$ws.Range("S7").Value = 4.4
$ws.Range("T7").Value = 1.95
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.44
$ws.Range("W7").Value = 1.6
$ws.Range("X7").Value = 11.5
$ws.Range("Y7").Value = 11.5
$ws.Range("Z7").Value = 21
$ws.Range("AB7").Value = 10
$ws.Range("AC7").Value = 7.4
$ws.Range("AD7").Value = 14.5
$ws.Range("AE7").Value = 42
$ws.Range("AF7").Value = 17
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 21
$ws.Range("AI7").Value = 65
$ws.Range("AJ7").Value = 42
$ws.Range("AL7").Value = 55
$ws.Range("AN7").Value = 34
$ws.Range("AO7").Value = 50
$ws.Range("J8").Value = 2.88
$ws.Range("K8").Value = 3.1
$ws.Range("N8").Value = 2.38
$ws.Range("S8").Value = 7
$ws.Range("Z8").Value = 14
$ws.Range("AE8").Value = 85
$ws.Range("G9").Value = 2.82
$ws.Range("H9").Value = 2.92
$ws.Range("I9").Value = 3.3
$ws.Range("J9").Value = 3.15
$ws.Range("K9").Value = 3.45
$ws.Range("L9").Value = 1.49
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 2.88
$ws.Range("S9").Value = 3.9
$ws.Range("U9").Value = 1.89
$ws.Range("V9").Value = 1.44
$ws.Range("W9").Value = 1.55
$ws.Range("AE9").Value = 55
$ws.Range("AI9").Value = 75
$ws.Range("AO9").Value = 60
$ws.Range("F10").Value = 1.23
$ws.Range("G10").Value = 1000
$ws.Range("H10").Value = 2.34
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 1.09
$ws.Range("K10").Value = 950
$ws.Range("O10").Value = 1.01
$ws.Range("Q10").Value = 1.01
$ws.Range("V10").Value = 1.34
$ws.Range("F11").Value = 2.18
$ws.Range("G11").Value = 2.38
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 3.35
$ws.Range("N11").Value = 2.68
$ws.Range("P11").Value = 1.64
$ws.Range("Q11").Value = 2.52
$ws.Range("T11").Value = 2.06
$ws.Range("U11").Value = 1.79
$ws.Range("W11").Value = 1.73
$ws.Range("AB11").Value = 7.6
$ws.Range("AC11").Value = 7.6
$ws.Range("J12").Value = 3.2
$ws.Range("L12").Value = 1.51
$ws.Range("O12").Value = 1.46
$ws.Range("Q12").Value = 2.24
$ws.Range("S12").Value = 4.5
$ws.Range("J13").Value = 3.55
$ws.Range("K13").Value = 3.95
$ws.Range("M13").Value = 1.07
$ws.Range("T13").Value = 1.87
$ws.Range("U13").Value = 1.93
$ws.Range("F14").Value = 1.47
$ws.Range("G14").Value = 1.65
$ws.Range("H14").Value = 6
$ws.Range("I14").Value = 9.4
$ws.Range("K14").Value = 5
$ws.Range("N14").Value = 3.9
$ws.Range("P14").Value = 2
$ws.Range("Q14").Value = 1.8
$ws.Range("R14").Value = 1.39
$ws.Range("W14").Value = 2.54
$ws.Range("X14").Value = 21
$ws.Range("Y14").Value = 980
$ws.Range("Z14").Value = 80
$ws.Range("AB14").Value = 10
$ws.Range("AC14").Value = 12.5
$ws.Range("AD14").Value = 36
$ws.Range("AE14").Value = 150
$ws.Range("AG14").Value = 12
$ws.Range("AL14").Value = 980
$ws.Range("AN14").Value = 9.800000000000001
$ws.Range("F15").Value = 2.6
$ws.Range("I15").Value = 3.45
$ws.Range("L15").Value = 1.62
$ws.Range("N15").Value = 2.44
$ws.Range("O15").Value = 1.59
$ws.Range("U15").Value = 1.74
$ws.Range("AJ15").Value = 980
$ws.Range("F16").Value = 1.72
$ws.Range("G16").Value = 1.76
$ws.Range("I16").Value = 6.8
$ws.Range("J16").Value = 3.65
$ws.Range("W16").Value = 2.32
